$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 189, shifting existing rows 189:273 down to 190:274.
$ws.Rows("189:189").Insert()

# Populate the newly inserted row 189 with the new record's data.
$ws.Range("A189").Value = 10
$ws.Range("B189").Value = "Vega Modelo de Temuco"
$ws.Range("C189").Value = "La Araucanía"
$ws.Range("D189").Value = 45016
$ws.Range("E189").Value = 9
$ws.Range("F189").Value = 100114007
$ws.Range("G189").Value = "Jengibre"
$ws.Range("H189").Value = "Sin especificar"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 80
$ws.Range("K189").Value = 25000
$ws.Range("L189").Value = 25000
$ws.Range("M189").Value = 25000
$ws.Range("N189").Value = "$/caja 13 kilos"
$ws.Range("O189").Value = "Perú"
$ws.Range("P189").Value = 1923
$ws.Range("Q189").Value = 13
$ws.Range("R189").Value = "Hortaliza"
